# update scripts with new tpm
# The NATMI LR-pair table (Fgf16 -> Fgfr2) was recomputed against a new
# TPM matrix. The ligand is always produced by "MuSCs" now (previously the
# sender/receiver pairing also included a "MuSCs -> *" block further down
# the sheet), so the 6 data rows collapse into 3 - one per target cluster
# (ECs, FAPs, MuSCs) - all carrying the freshly recomputed metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#          Ligand-expressing cells, Ligand detection rate,
#          Ligand average expression value, Ligand total expression value,
#          Ligand derived specificity (avg), Ligand derived specificity (total),
#          Receptor-expressing cells, Receptor detection rate,
#          Receptor average expression value, Receptor total expression value,
#          Receptor derived specificity (avg), Receptor derived specificity (total),
#          Edge average expression weight, Edge total expression weight,
#          Edge average expression derived specificity, Edge total expression derived specificity
$data = @(
    @("MuSCs","Fgf16","Fgfr2","ECs",  3,1,0.2782933333333333,0.83488,1,1,3,1,0.106124,0.318372,0.08094716512538251,0.08094716512538253,0.02953360170666667,0.26580241536,0.08094716512538251,0.08094716512538253),
    @("MuSCs","Fgf16","Fgfr2","FAPs", 3,1,0.2782933333333333,0.83488,1,1,3,1,1.092289666666667,3.276869,0.8331551016962769,0.833155101696277,0.3039769323022222,2.73579239072,0.8331551016962769,0.833155101696277),
    @("MuSCs","Fgf16","Fgfr2","MuSCs",3,1,0.2782933333333333,0.83488,1,1,3,1,0.1126143333333333,0.337843,0.08589773317834044,0.08589773317834046,0.03133981820444445,0.28205836384,0.08589773317834044,0.08589773317834046)
)

# Wipe out the old 6-row table body (rows 2..7), leave the header row intact.
$ws.Range("A2:T7").Clear()

# Write the new 3-row table body back in starting at row 2.
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
